$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price observation was added to the data set. It belongs
# right before the existing row 354 (by date order within the sheet), so
# insert a whole new row there, which pushes the former rows 354-413 down
# to become rows 355-414 (all their data travels with them unchanged).
$ws.Rows.Item(354).Insert()

# Populate the newly inserted row 354 with the new observation's data.
$ws.Cells.Item(354, 1).Value = 8
$ws.Cells.Item(354, 2).Value = 'Terminal La Palmera de La Serena'
$ws.Cells.Item(354, 3).Value = 'Coquimbo'
$ws.Cells.Item(354, 4).Value = 44951
$ws.Cells.Item(354, 5).Value = 4
$ws.Cells.Item(354, 6).Value = 100112003
$ws.Cells.Item(354, 7).Value = 'Ajo'
$ws.Cells.Item(354, 8).Value = 'Chino'
$ws.Cells.Item(354, 9).Value = 'Primera'
$ws.Cells.Item(354, 10).Value = 480
$ws.Cells.Item(354, 11).Value = 16000
$ws.Cells.Item(354, 12).Value = 17000
$ws.Cells.Item(354, 13).Value = 16500
$ws.Cells.Item(354, 14).Value = '$/caja 10 kilos'
$ws.Cells.Item(354, 15).Value = 'China'
$ws.Cells.Item(354, 16).Value = 1650
$ws.Cells.Item(354, 17).Value = 10
$ws.Cells.Item(354, 18).Value = 'Hortaliza'
